# Mise à jour de l'application
# Adds a new training-session column (BG) dated 2025-10-07 to the attendance
# sheet, records each player's status for that session, and appends a
# totals row (row 30) that counts the "P" (présent) entries for BG.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New training date header (row 1) ---------------------------------
# Write the value first (so the recalculation graph picks the cell up),
# then copy BF1's style (date format, centered) onto BG1, and finally
# restore the value (Copy also copies BF1's value, which we overwrite).
$ws.Range("BG1").Value = 45937
$ws.Range("BF1").Copy($ws.Range("BG1"))
$ws.Range("BG1").Value = 45937

# --- Per-player attendance status for the new session (column BG) -----
# Row 12 (Yanis Berrached) has no data past column AX, so it is excluded.
$statuses = @{
  2="P"; 3="R"; 4="P"; 5="P"; 6="P"; 7="P"; 8="B"; 9="P"; 10="P"; 11="P";
  13="B"; 14="P"; 15="P"; 16="B"; 17="P"; 18="P"; 19="M"; 20="P"; 21="M";
  22="P"; 23="B"; 24="P"; 25="P"; 26="P"; 27="P"; 28="P"; 29="P"
}

foreach ($r in $statuses.Keys) {
  $src = $ws.Range("BF" + $r)
  $dst = $ws.Range("BG" + $r)
  $val = $statuses[$r]
  $dst.Value = $val
  $src.Copy($dst)
  $dst.Value = $val
}

# --- New totals row (row 30): count of "P" for the new column ---------
$ws.Range("BG30").Formula = "=COUNTIF(BG5:BG29,""P"")"

# --- View state: keep the frozen first column and reselect the bottom
# of the newly-added column, matching the author's saved selection. ----
$ws.Range("BG29").Select()
